$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Team Info": turn the two-season summary into a three-season one and
# move the "Number of Seasons" counter from K1:K2 up to A1:A2.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Team Info")

# Push the two existing season blocks down by two rows, freeing rows 1:2 for
# the season counter.
$ws.Range("A1:A2").EntireRow.Insert()

# Move "Number of Seasons" / 3 from K3:K4 (the old K1:K2, shifted down by the
# insert above) up to A1:A2 -- this also clears out column K.
$ws.Range("K3:K4").Cut($ws.Range("A1:A2"))

# Make room at the bottom of the sheet for a brand-new Season 3 block: one
# blank separator row, a "Season No./Player Count" header row, a value row,
# a "Player Names" label row and a names row.
$ws.Range("A12:A16").EntireRow.Insert()

# Header row ("Season No." / "Player Count"), formatted like the other two
# season headers.
$ws.Range("A8:B8").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A13").Value = "Season No."
$ws.Range("B13").Value = "Player Count"
$ws.Range("A14").Value = 3
$ws.Range("B14").Value = 15

# "Player Names" label row.
$ws.Range("A10").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A15").Value = "Player Names"

# Player-name row (same ten names as the other two seasons).
$ws.Range("A11:J11").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A16").Value = "Brandon Chan"
$ws.Range("B16").Value = "Callum Ashton"
$ws.Range("C16").Value = "Daniel Park"
$ws.Range("D16").Value = "Deidre Truong"
$ws.Range("E16").Value = "Edward Kang"
$ws.Range("F16").Value = "Kevin Ma"
$ws.Range("G16").Value = "Kevin Tang"
$ws.Range("H16").Value = "Lachlan Denham"
$ws.Range("I16").Value = "Mimi Chen"
$ws.Range("J16").Value = "Will Ouyang"

# Scroll the sheet so the new Season 3 block is in view, and select A18:B22
# (the next free rows below it) -- where a user would type a fourth season.
$ws.Application.Goto($ws.Range("A16"))
$ws.Range("A18:B22").Select()

# ---------------------------------------------------------------------------
# Sheet "Season 1": scroll down near the bottom of the entered stats.
# ---------------------------------------------------------------------------
$s1 = $wb.Worksheets.Item("Season 1")
$s1.Application.Goto($s1.Range("A100"))
$s1.Range("A3:A12").Select()

# ---------------------------------------------------------------------------
# Sheet "Season 3": move the selection to B27, ready for new data entry.
# ---------------------------------------------------------------------------
$s3 = $wb.Worksheets.Item("Season 3")
$s3.Range("B27").Select()

$ws.Activate()
